$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Rename header row strings: _old -> _FV2404, _new -> _FV2410
$oldNames = @("Segmentname_old","Segmentgruppe_old","Segment_old","Datenelement_old","Segment ID_old","Code_old","Qualifier_old","Beschreibung_old","Bedingungsausdruck_old","Bedingung_old")
$newFV2404 = @("Segmentname_FV2404","Segmentgruppe_FV2404","Segment_FV2404","Datenelement_FV2404","Segment ID_FV2404","Code_FV2404","Qualifier_FV2404","Beschreibung_FV2404","Bedingungsausdruck_FV2404","Bedingung_FV2404")

$newNames = @("Segmentname_new","Segmentgruppe_new","Segment_new","Datenelement_new","Segment ID_new","Code_new","Qualifier_new","Beschreibung_new","Bedingungsausdruck_new","Bedingung_new")
$newFV2410 = @("Segmentname_FV2410","Segmentgruppe_FV2410","Segment_FV2410","Datenelement_FV2410","Segment ID_FV2410","Code_FV2410","Qualifier_FV2410","Beschreibung_FV2410","Bedingungsausdruck_FV2410","Bedingung_FV2410")

for ($i = 0; $i -lt 10; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $newFV2404[$i]
}
# column 11 is "diff" (unchanged)
for ($i = 0; $i -lt 10; $i++) {
    $ws.Cells.Item(1, $i + 12).Value = $newFV2410[$i]
}

# 2) Add an Excel Table (ListObject) over A1:U69 with header row
$tbl = $ws.ListObjects.Add(1, $ws.Range("A1:U69"), 0, 1)
$tbl.Name = "Table1"

# 3) Freeze the top row (header)
$ws.Activate()
$excel.ActiveWindow.SplitRow = 1
$excel.ActiveWindow.SplitColumn = 0
$excel.ActiveWindow.FreezePanes = $true
